$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G15").Value = '2023/2024, 2025/2026'
$ws.Range("G19").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G21").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G22").Value = '2025/2026, Eman_mohamed@med.asu.edu.eg'
$ws.Range("G37").Value = '2023/2024, 2025/2026'
$ws.Range("G41").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G43").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G44").Value = '2025/2026, Eman_mohamed@med.asu.edu.eg'
$ws.Range("G60").Value = '2025/2026, 2026/2027'
$ws.Range("G63").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G64").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G65").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G66").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G82").Value = '2025/2026, 2026/2027'
$ws.Range("G85").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G86").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G87").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G88").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G89").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G104").Value = '2025/2026, 2026/2027'
$ws.Range("G106").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G107").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G110").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G111").Value = 'youstina.magdy@med.asu.edu.eg, 2025/2026, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg'
$ws.Range("G126").Value = '2025/2026, 2026/2027'
$ws.Range("G128").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G129").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G132").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G133").Value = 'youstina.magdy@med.asu.edu.eg, 2025/2026, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg'
$ws.Range("G150").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G153").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G155").Value = 'youstina.magdy@med.asu.edu.eg, 2025/2026, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg'
$ws.Range("G172").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G175").Value = 'neveen.nashaat@med.asu.edu.eg, 2025/2026'
$ws.Range("G177").Value = 'youstina.magdy@med.asu.edu.eg, 2025/2026, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg'
